# Apply "Add data for 2022-05-30" update:
# - Rename sheet/title from "...May 21" to "...May 22"
# - Bump the May-2022 header label
# - Update carjacking counts for several neighborhoods' May columns

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab (drives the <sheet name="..."> entry).
$ws.Name = "Through 2022-05-22"

# Update the column header text for the "May 2022" running total.
$ws.Range("B1").Value = "May 2022 (through May 22)"

# --- Row 2: Englewood ---
$ws.Range("B2").Value = 8   # May 2022 (through May 22): 7 -> 8
$ws.Range("G2").Value = 4   # May 2021: 3 -> 4

# --- Row 3: Austin ---
$ws.Range("G3").Value = 7   # May 2021: 6 -> 7

# --- Row 5: Garfield Park ---
$ws.Range("B5").Value = 2   # May 2022 (through May 22): 1 -> 2
$ws.Range("Q5").Value = 2   # May 2019: new value 2

# --- Row 8: South Shore ---
$ws.Range("Q8").Value = 2   # May 2019: 1 -> 2

# --- Row 20: Woodlawn ---
$ws.Range("V20").Value = 2  # May 2018: 1 -> 2

# --- Row 21: Chatham ---
$ws.Range("AK21").Value = 1 # May 2015: new value 1

# --- Row 23: Grand Crossing ---
$ws.Range("B23").Value = 4  # May 2022 (through May 22): 3 -> 4
$ws.Range("AA23").Value = 2 # May 2017: 1 -> 2

# --- Row 31: Uptown ---
$ws.Range("AA31").Value = 1 # May 2017: new value 1

# --- Row 45: Logan Square ---
$ws.Range("B45").Value = 2  # May 2022 (through May 22): 1 -> 2

# --- Row 90: Ukrainian Village ---
$ws.Range("AA90").Value = 1 # May 2017: new value 1
